$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 to have spaces in the NIF number
$ws.Range("A2").Value = "309 889 669"

# Add new row 3 with data
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "313424640"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "01"
$ws.Range("C3").Value = 180
